$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New headers: "Area" (col G) / "Atotal" (col H)
$ws.Range("G1").Value = "Area"
$ws.Range("H1").Value = "Atotal"

# G2: area of the first segment, measured from 0
$ws.Range("G2").Formula = "=(D2-0)*B2/100"

# H2: running total of the area column
$ws.Range("H2").Formula = "=SUM(G2:G11)"

# G3: area of the second segment (first cell using the D-Dprev pattern)
$ws.Range("G3").Formula = "=(D3-D2)*B3/100"

# G4:G15 continue the same relative "(D-Dprev)*B/100" pattern
for ($r = 4; $r -le 15; $r++) {
    $prev = $r - 1
    $ws.Range("G$r").Formula = "=(D$r-D$prev)*B$r/100"
}

$ws.Calculate()

# Match the saved selection/scroll state from the edit
$ws.Range("H2").Select()
